$d = $word.ActiveDocument

# wdReplace / wdFind constants
$wdReplaceOne = 1
$wdFindContinue = 1

# --- Figure S1 (B) caption: hypergraph conversion sentence ---
$old1 = "The metabolic network reconstruction is converted to a hypergraph, in which metabolites are represented as nodes and reactions as hyperedges. In this representation, an edge can connect more than two nodes. For clarity, protons are not shown."
$new1 = "The metabolic network reconstruction represents metabolism as a hypergraph, in which metabolites are represented as nodes and reactions as hyperedges. In this representation, an edge can connect more than two nodes. For example, a single hyperedge (denoted by a heavy black line) connects the metabolites glucose and ATP to glucose-6P, ADP, and Pi. For clarity, protons are not shown."

$ok1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new1, $wdReplaceOne)
Write-Output "Replace 1 (hypergraph sentence): $ok1"

# --- Figure S1 (C) caption: metabolic network graph sentence ---
$old2 = "The hypergraph is converted to a metabolic network graph, in which an edge can connect only two nodes. In this representation, a reaction is represented by a set of edges connecting all substrates to all products. The dotted line surrounds the currency metabolites."
$new2 = "However, the algorithm used by the seed set framework requires metabolism to be represented as a metabolic network graph, in which an edge can connect only two nodes. In this representation, a reaction is represented by a set of edges connecting all substrates to all products. For example, the heavy hyperedge in (B) is now denoted by six separate edges connecting glucose to ADP, glucose to Pi, glucose to glucose-6P, ATP to ADP, ATP to Pi, and ATP to glucose-6P (again denoted by heavy black lines). Of these, only one (glucose to glucose-6P) is biologically meaningful. The dotted line surrounds the currency metabolites."

$ok2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new2, $wdReplaceOne)
Write-Output "Replace 2 (metabolic network graph sentence): $ok2"

# --- Figure S1 (D) caption: pruning sentence ---
$old3 = "The metabolic network graph is then pruned, a process which removes all currency metabolites and any edges in which those metabolites participate. Representation of glycolysis after pruning. The images in (B) and (C) are modified from"
$new3 = "The metabolic network graph is then pruned, a process which removes all currency metabolites and any edges in which those metabolites participate. Of the six heavy edges in (C), only the biologically meaningful one is retained, connecting glucose to glucose-6P (again denoted by a heavy black line). The images in (B) and (C) are modified from"

$ok3 = $d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new3, $wdReplaceOne)
Write-Output "Replace 3 (pruning sentence): $ok3"
